$d = $word.ActiveDocument

function Get-ParaByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pp = $d.Paragraphs.Item($i)
        $t = $pp.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $pp
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1) Delete the paragraph "*Chức năng đăng kí hợp đồng kí túc xá"
# ------------------------------------------------------------------
$star = Get-ParaByText "*Chức năng đăng kí hợp đồng kí túc xá"
if ($star -ne $null) {
    $star.Range.Delete()
}

# ------------------------------------------------------------------
# 2) Bump "I.Mô tả chức năng " heading from 28 -> 32 half-points
#    (i.e. 14pt -> 16pt), both paragraph-mark formatting and the run.
# ------------------------------------------------------------------
$heading = Get-ParaByText "I.Mô tả chức năng "
$heading.Range.Font.Size = 16
$heading.Range.Font.SizeBi = 16

# ------------------------------------------------------------------
# 3) Insert three new paragraphs right after the heading:
#      "1." (bold) "Đăng nhập hệ thống" (regular)
#      <tab> + relocated _GoBack bookmark
#      "2." (bold) "Đăng kí hợp Đồng" (regular)
#    Done by inserting a throw-away empty paragraph right after the
#    heading and then replacing that whole (start..end) empty range
#    with a raw-XML fragment describing the three paragraphs exactly,
#    so the paragraph-mark / run formatting matches precisely (no
#    forced-bold inheritance, no explicit b val="0").
# ------------------------------------------------------------------
$heading.Range.InsertParagraphAfter()
$heading = Get-ParaByText "I.Mô tả chức năng "
$spacer = $heading.Next()

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$p1 = "<w:p $ns><w:pPr><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr>" +
      "<w:r><w:rPr><w:b/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>1.</w:t></w:r>" +
      "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>Đăng nhập hệ thống</w:t></w:r></w:p>"
$p2 = "<w:p $ns><w:pPr><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr>" +
      "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:tab/></w:r>" +
      "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$p3 = "<w:p $ns><w:pPr><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr>" +
      "<w:r><w:rPr><w:b/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>2.</w:t></w:r>" +
      "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>Đăng kí hợp Đồng</w:t></w:r></w:p>"

$body = $p1 + $p2 + $p3
$fullXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           "<w:body>$body</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$rng = $d.Range($spacer.Range.Start, $spacer.Range.End)
$rng.InsertXML($fullXml)

# ------------------------------------------------------------------
# 4) Remove the old _GoBack bookmark that used to sit after
#    "II.Công nghệ và công cụ" (the bookmark now lives in the new
#    paragraph we just inserted above).
# ------------------------------------------------------------------
$techHeading = Get-ParaByText "II.Công nghệ và công cụ"
if ($d.Bookmarks.Count -gt 1) {
    for ($i = $d.Bookmarks.Count; $i -ge 1; $i--) {
        $bm = $d.Bookmarks.Item($i)
        if ($bm.Name -eq "_GoBack" -and $bm.Start -ge $techHeading.Range.Start) {
            $bm.Delete()
        }
    }
}

# ------------------------------------------------------------------
# 5) Add a <w:lastRenderedPageBreak/> right before the bold "-" run
#    that starts the "-công cụ : ..." paragraph.
# ------------------------------------------------------------------
$toolsPara = Get-ParaByText "-công cụ : netbean  và  SQL Sever"
if ($toolsPara -ne $null) {
    $firstRun = $d.Range($toolsPara.Range.Start, $toolsPara.Range.Start + 1)
    $runXml = $firstRun.WordOpenXML
}
Write-Output "done"
